$wb = $excel.ActiveWorkbook

# Companies sheet: move the selection from A4 to A3
$wsCompanies = $wb.Worksheets.Item("Companies")
$wsCompanies.Activate()
$wsCompanies.Range("A3").Select() | Out-Null

# Contacts sheet: add the new contact name and update the selection,
# then leave it as the active sheet/tab
$wsContacts = $wb.Worksheets.Item("Contacts")
$wsContacts.Activate()
$wsContacts.Range("B2").Value = "Adam Cole"
$wsContacts.Range("A2:B2").Select() | Out-Null
